# issue #5: property land done
#
# Normalizes whitespace/punctuation noise in the "土地" (land), "建物"
# (building) and "債務" (debt) sheets, and appends the canonical English
# metadata columns (property_category/category/date/legislator_name/
# legislator_id/source_file/index) to the "土地" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write $text into $range as a literal text value even when it
# looks numeric/date-like (Excel would otherwise coerce "6000000" or
# "2012-04-30" into a real number/date). Round-tripping through a TEXT()
# formula and then pasting-special as values keeps the cell's existing
# style untouched (no NumberFormat/quote-prefix residue).
# ---------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text -replace '"', '""'
    $range.Formula = '=TEXT("' + $escaped + '","@")'
    $range.Copy()
    $range.PasteSpecial(-4163) | Out-Null
}

# =======================================================================
# Sheet "土地" (land)
# =======================================================================
$landSheet = $wb.Worksheets.Item("土地")

# -- normalize existing header/data text (strip inserted spaces/dashes/commas)
$landSheet.Range("B1").Value = "name"
$landSheet.Range("C1").Value = "area"
$landSheet.Range("D1").Value = "share_portion"
$landSheet.Range("E1").Value = "owner"
$landSheet.Range("F1").Value = "register_date"
$landSheet.Range("G1").Value = "register_reason"
$landSheet.Range("H1").Value = "acquire_value"

$landSheet.Range("B2").Value = "臺北市大安區龍泉段一小段02930000地號"
$landSheet.Range("D2").Value = "100000分之16216"
$landSheet.Range("F2").Value = "92年12月25日"
$landSheet.Range("G2").Value = "033貝賣"
$landSheet.Range("H2").Value = "25000000(土地建物與車位合併價）"

# -- append the new metadata columns I:O, copying style from the existing
#    header (B1, bold+bordered) / data (B2, plain) cells so the new cells
#    match the sheet's existing look.
$landSheet.Range("B1").Copy($landSheet.Range("I1"))
$landSheet.Range("I1").Value = "property_category"
$landSheet.Range("B1").Copy($landSheet.Range("J1"))
$landSheet.Range("J1").Value = "category"
$landSheet.Range("B1").Copy($landSheet.Range("K1"))
$landSheet.Range("K1").Value = "date"
$landSheet.Range("B1").Copy($landSheet.Range("L1"))
$landSheet.Range("L1").Value = "legislator_name"
$landSheet.Range("B1").Copy($landSheet.Range("M1"))
$landSheet.Range("M1").Value = "legislator_id"
$landSheet.Range("B1").Copy($landSheet.Range("N1"))
$landSheet.Range("N1").Value = "source_file"
$landSheet.Range("B1").Copy($landSheet.Range("O1"))
$landSheet.Range("O1").Value = "index"

$landSheet.Range("B2").Copy($landSheet.Range("I2"))
$landSheet.Range("I2").Value = "land"
$landSheet.Range("B2").Copy($landSheet.Range("J2"))
$landSheet.Range("J2").Value = "normal"
$landSheet.Range("B2").Copy($landSheet.Range("K2"))
Set-TextValue $landSheet.Range("K2") "2012-04-30"
$landSheet.Range("B2").Copy($landSheet.Range("L2"))
$landSheet.Range("L2").Value = "高金素梅"
$landSheet.Range("B2").Copy($landSheet.Range("M2"))
$landSheet.Range("M2").Value = 926
$landSheet.Range("B2").Copy($landSheet.Range("N2"))
$landSheet.Range("N2").Value = "tmpb18e1"
$landSheet.Range("B2").Copy($landSheet.Range("O2"))
$landSheet.Range("O2").Value = 14

# =======================================================================
# Sheet "建物" (building) - normalize existing text only
# =======================================================================
$buildingSheet = $wb.Worksheets.Item("建物")
$buildingSheet.Range("B2").Value = "臺北市大安區龍泉段一小段05819000建號"
$buildingSheet.Range("F2").Value = "92年12月25日"
$buildingSheet.Range("H2").Value = "25000000(土地建物與車位合併價）"

# =======================================================================
# Sheet "存款" (deposit) - unchanged (no cell edits required)
# =======================================================================

# =======================================================================
# Sheet "債務" (debt) - normalize existing text only
# =======================================================================
$debtSheet = $wb.Worksheets.Item("債務")
$debtSheet.Range("D2").Value = "陳麗卿新北市泰山區明志路"
Set-TextValue $debtSheet.Range("E2") "6000000"
$debtSheet.Range("F2").Value = "96年02月06日"

$debtSheet.Range("D3").Value = "石旭松新北市泰山區明志路"
Set-TextValue $debtSheet.Range("E3") "4000000"
$debtSheet.Range("F3").Value = "96年02月06日"
